{"js": "// Add a new paragraph style \"Footnote Block Text\", based on (and followed\n// by) \"Footnote Text\", mirroring the existing \"Block Text\" style but for\n// block quotes that live inside footnotes/endnotes.\ncontext.document.addStyle(\"Footnote Block Text\", Word.StyleType.paragraph);\nawait context.sync();\n\n// Re-fetch the style from the styles collection by name; writes only stick\n// reliably once the style is addressed through the collection (rather than\n// the anchor handed back directly by addStyle()).\nconst style = context.document.getStyles().getByName(\"Footnote Block Text\");\n\nstyle.baseStyle = \"FootnoteText\";\nstyle.nextParagraphStyle = \"FootnoteText\";\nstyle.priority = 9;\nstyle.unhideWhenUsed = true;\nstyle.quickStyle = true;\n\n// Match \"Block Text\"'s paragraph formatting: spacing before/after 100\n// twips (5pt) and left/right indents of 480 twips (24pt), no first-line\n// indent.\nstyle.paragraphFormat.spaceBefore = 5;\nstyle.paragraphFormat.spaceAfter = 5;\nstyle.paragraphFormat.firstLineIndent = 0;\nstyle.paragraphFormat.leftIndent = 24;\nstyle.paragraphFormat.rightIndent = 24;\n\nawait context.sync();\n", "ps1": "# Add a new paragraph style \"Footnote Block Text\", based on (and followed\n# by) \"Footnote Text\", mirroring the existing \"Block Text\" style but for\n# block quotes that live inside footnotes/endnotes.\n$d = $word.ActiveDocument\n\n$style = $d.Styles.Add(\"Footnote Block Text\", \"wdStyleTypeParagraph\")\n$footnoteText = $d.Styles(\"FootnoteText\")\n\n$style.BaseStyle = $footnoteText\n$style.NextParagraphStyle = $footnoteText\n$style.Priority = 9\n$style.UnhideWhenUsed = $true\n$style.QuickStyle = $true\n\n# Match \"Block Text\"'s paragraph formatting: spacing before/after 100\n# twips (5pt) and left/right indents of 480 twips (24pt), no first-line\n# indent.\n$style.ParagraphFormat.SpaceBefore = 5\n$style.ParagraphFormat.SpaceAfter = 5\n$style.ParagraphFormat.FirstLineIndent = 0\n$style.ParagraphFormat.LeftIndent = 24\n$style.ParagraphFormat.RightIndent = 24\n"}
